# OCT 25 update, removed phone number
# ------------------------------------------------------------------
# Adds a new "Active" status column (F) to the collaborations table,
# corrects the "Legend Biotech" collaboration row (row 18), and
# appends two new collaboration rows (Luara Gleeson / Kathy Gately).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A style-"5" cell (plain, vertically centred, no wrap) used as a format
# donor for the new column-F cells that need to match the rest of the
# table's "year"/"where" styling.
$styleDonor5 = $ws.Range("C5")

# --- New column F header --------------------------------------------------
$ws.Range("F1").Value = "Active"

# --- Populate column F for existing data rows (style donor first, then
#     value, so the copy doesn't clobber the text we are about to write) --
$styleDonor5.Copy($ws.Range("F5"))
$styleDonor5.Copy($ws.Range("F7"))
$styleDonor5.Copy($ws.Range("F8"))
$styleDonor5.Copy($ws.Range("F9"))
$styleDonor5.Copy($ws.Range("F10"))
$styleDonor5.Copy($ws.Range("F11"))
$styleDonor5.Copy($ws.Range("F12"))
$styleDonor5.Copy($ws.Range("F13"))
$styleDonor5.Copy($ws.Range("F14"))
$styleDonor5.Copy($ws.Range("F15"))
$styleDonor5.Copy($ws.Range("F16"))
$styleDonor5.Copy($ws.Range("F17"))
$styleDonor5.Copy($ws.Range("F18"))

$ws.Range("F4").Value  = "Y"
$ws.Range("F5").Value  = "Y"
$ws.Range("F6").Value  = "Y"
$ws.Range("F7").Value  = "N"
$ws.Range("F8").Value  = "Y"
$ws.Range("F9").Value  = "Y"
$ws.Range("F10").Value = "Y"
$ws.Range("F11").Value = "N"
$ws.Range("F12").Value = "Y"
$ws.Range("F13").Value = "Y"
$ws.Range("F14").Value = "Y"
$ws.Range("F15").Value = "Y"
$ws.Range("F16").Value = "Y"
$ws.Range("F17").Value = "Y"
$ws.Range("F18").Value = "Y"

# --- Row 18: Legend Biotech collaboration - corrected details ------------
$ws.Range("A18").Value = "Legend Biotech, Tony Mc Elligott & Nina Orfali"
$ws.Range("B18").Value = "Blood Cancer, Title withheld"
$ws.Range("D18").Value = "TCD"
$ws.Range("E18").Value = "Title withheald"

# --- Row 19: new collaboration (Luara Gleeson) ----------------------------
# Seed formatting from row 18 (same column styles apply to row 19).
$ws.Range("A18").Copy($ws.Range("A19"))
$ws.Range("C18").Copy($ws.Range("D19"))
$ws.Range("E18").Copy($ws.Range("E19"))
$ws.Range("C18").Copy($ws.Range("F19"))

$ws.Range("A19").Value = "Luara Gleeson"
$ws.Range("C19").Value = "2024-Present"
$ws.Range("D19").Value = "TCD"
$ws.Range("E19").Value = "Human Pleural biology"
$ws.Range("F19").Value = "Y"
$ws.Rows.Item(19).RowHeight = 16

# --- Row 20: new collaboration (Kathy Gately) -----------------------------
$ws.Range("A18").Copy($ws.Range("A20"))
$ws.Range("C18").Copy($ws.Range("C20"))
$ws.Range("C18").Copy($ws.Range("D20"))
$ws.Range("E18").Copy($ws.Range("E20"))
$ws.Range("C18").Copy($ws.Range("F20"))

$ws.Range("A20").Value = "Kathy Gately "
$ws.Range("C20").Value = "2025-Present"
$ws.Range("D20").Value = "TCD, St James' Hospital"
$ws.Range("E20").Value = "Malignang Pleural Effussion"
$ws.Range("F20").Value = "Y"
$ws.Rows.Item(20).RowHeight = 16

# --- Selection / view state: column F highlighted, as in the authored file
$ws.Application.Goto($ws.Range("F1:F19"), $true)
